# user_profiles.xlsx - "validation and excel format"
#
# Reformat the column headers from snake_case field names into
# human-friendly Title Case labels, and restore the sheet's last
# on-screen selection/formatting state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Excel format: snake_case headers -> Title Case -------------------
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Middle Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "Mobile"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Password"

# --- Keep the bottom-most formatted row in sync with the sheet's ------
# default row formatting (matches the sheet's existing last row).
$ws.Rows.Item(1048575).RowHeight = $ws.Rows.Item(1048576).RowHeight

# --- Restore the active cell/selection saved with the workbook --------
$ws.Range("C14").Select()
